$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix of the two same subjects - correct the duplicated/shifted
# "Tecnologia dos Materiais" entries in column D (Wednesday).
$ws.Range("D2").Value = "MEC-1A-Tecnologia dos Materiais"
$ws.Range("D3").Value = "MEC-1A-Tecnologia dos Materiais"
$ws.Range("D4").Value = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("D6").Value = "-"
